$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A4").Value = 43819
$ws.Range("B4").Value = 7949
$ws.Range("C4").Value = 46315
$ws.Range("D4").Value = 3.651621092140373
$ws.Range("E4").Value = 0.6624621482598401
$ws.Range("F4").Value = 3.859650691720534
$ws.Range("G4").Value = 0.2532260257314133
$ws.Range("H4").Value = 0.07494766662997333
$ws.Range("I4").Value = 0.22811521316864
$ws.Range("J4").Value = 0.2194322041714432
$ws.Range("K4").Value = 0.1494937248756479
$ws.Range("L4").Value = 0.2189282773549824
$ws.Range("M4").Value = 51433
$ws.Range("N4").Value = 7570
$ws.Range("O4").Value = 122848
$ws.Range("P4").Value = 4.485457804328482
$ws.Range("Q4").Value = 0.6601804948741952
$ws.Range("R4").Value = 10.71351860878884
$ws.Range("S4").Value = 0.3360088773870139
$ws.Range("T4").Value = 0.07009107148085582
$ws.Range("U4").Value = 0.5365403292529116
$ws.Range("V4").Value = 0.3030822041714432
$ws.Range("W4").Value = 0.1600382520029951
$ws.Range("X4").Value = 0.3850826323743871
